# Apply the edit described by the diff:
#  - Header labels on both sheets change suffix from "CL" to "CI"
#    (LL95CL -> LL95CI, UL95CL -> UL95CI, LL998CL -> LL998CI, UL998CL -> UL998CI)
#  - Active sheet / tab changes from "testdata_Rate_100000" (sheet 2)
#    back to "testdata_Rate_100" (sheet 1)
#  - Selection (active cell) on each sheet is updated

$wb = $excel.ActiveWorkbook
$s1 = $wb.Sheets.Item(1)   # testdata_Rate_100
$s2 = $wb.Sheets.Item(2)   # testdata_Rate_100000

# Update the confidence-limit header labels (CL -> CI) on sheet 1
$s1.Range("E1").Value = "LL95CI"
$s1.Range("F1").Value = "UL95CI"
$s1.Range("G1").Value = "LL998CI"
$s1.Range("H1").Value = "UL998CI"

# Update the confidence-limit header labels (CL -> CI) on sheet 2
$s2.Range("E1").Value = "LL95CI"
$s2.Range("F1").Value = "UL95CI"
$s2.Range("G1").Value = "LL998CI"
$s2.Range("H1").Value = "UL998CI"

# Move the selection on sheet 2 (no longer the active tab) before leaving it
$s2.Activate()
$s2.Range("L10").Select()

# Sheet 1 becomes the active tab again, with its own updated selection
$s1.Activate()
$s1.Range("F12").Select()
